# Update "想去人数" (column F) counts across the four worksheets to reflect
# the latest scrape of the site data (gh-pages output regeneration).

$wb = $excel.ActiveWorkbook

# 展览 (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value  = 561
$ws1.Range("F8").Value  = 1209
$ws1.Range("F9").Value  = 974
$ws1.Range("F10").Value = 443
$ws1.Range("F12").Value = 688
$ws1.Range("F15").Value = 744
$ws1.Range("F16").Value = 19
$ws1.Range("F18").Value = 1316
$ws1.Range("F19").Value = 20
$ws1.Range("F24").Value = 299
$ws1.Range("F25").Value = 420
$ws1.Range("F27").Value = 69
$ws1.Range("F29").Value = 1
$ws1.Range("F32").Value = 260
$ws1.Range("F33").Value = 143
$ws1.Range("F34").Value = 53

# 演出 (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F9").Value  = 653
$ws2.Range("F10").Value = 368
$ws2.Range("F20").Value = 577
$ws2.Range("F23").Value = 422
$ws2.Range("F25").Value = 7
$ws2.Range("F26").Value = 180

# 本地生活 (Local life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F6").Value  = 2134
$ws3.Range("F7").Value  = 838
$ws3.Range("F8").Value  = 787
$ws3.Range("F9").Value  = 1
$ws3.Range("F11").Value = 795
$ws3.Range("F12").Value = 120

# 全部类型 (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value  = 2134
$ws4.Range("F8").Value  = 838
$ws4.Range("F9").Value  = 787
$ws4.Range("F13").Value = 561
$ws4.Range("F16").Value = 1209
$ws4.Range("F17").Value = 974
$ws4.Range("F18").Value = 795
$ws4.Range("F19").Value = 443
$ws4.Range("F21").Value = 120
$ws4.Range("F23").Value = 653
$ws4.Range("F24").Value = 744
$ws4.Range("F26").Value = 20
$ws4.Range("F31").Value = 299
$ws4.Range("F32").Value = 420
$ws4.Range("F33").Value = 69
$ws4.Range("F35").Value = 1
$ws4.Range("F37").Value = 260
$ws4.Range("F38").Value = 143
$ws4.Range("F39").Value = 53
$ws4.Range("F44").Value = 7
$ws4.Range("F45").Value = 180
